# Generate Report for Handoff
# Updates the localization-status report:
#   - Status "In Translation" -> "Ready for handoff"
#   - Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps bumped
#   - Widen the "Status" column(s) to fit the new, longer status text

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 06:52:06"
$wsOverview.Range("E1:F1").ColumnWidth = 16.3

# ---- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 06:51:58"
$wsZhCn.Range("C1").ColumnWidth = 16.3

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 06:52:06"
$wsDeDe.Range("C1").ColumnWidth = 16.3
